# Fix problem of multiple worksheets being instantiated:
# remove the duplicate "t2" (Televisor) row and the duplicate "l1" (Lampada) row,
# and update the remaining sensor/automation values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra "Televisor" row (originally row 6) first so row numbers
# above it are unaffected, then delete the extra "Lampada" row (originally row 3).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(3).Delete()

# Update remaining values to reflect the corrected data.
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = $false

$ws.Range("C2").Value = 18

$ws.Range("C3").Value = 50
